# Refresh cryptos list: update Price / Volume(1h) columns and fix
# the WEMIXToken / EthereumClassic row ordering (rows 32-33).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.548.30'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '2.224.74'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'111.89"
$ws.Range('E5').Value = '  -1.54%  '
$ws.Range('D6').Value = "'293.94"
$ws.Range('E6').Value = '  +10.15%  '
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('D9').Value = "'0.597"
$ws.Range('E9').Value = '  -1.46%  '
$ws.Range('D10').Value = "'43.33"
$ws.Range('E10').Value = '  -6.31%  '
$ws.Range('D11').Value = "'0.0919"
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').Value = "'54.20"
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('D13').Value = "'8.68"
$ws.Range('E13').Value = '  -4.85%  '
$ws.Range('D14').Value = "'1.05"
$ws.Range('E14').Value = '  +19.64%  '
$ws.Range('E15').Value = '  -1.36%  '
$ws.Range('E16').Value = '  -2.33%  '
$ws.Range('D17').Value = '2.560.77'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').Value = '2.217.79'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('D19').Value = '42.383.70'
$ws.Range('E19').Value = '  -1.39%  '
$ws.Range('D20').Value = "'7.18"
$ws.Range('E20').Value = '  +6.76%  '
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').Value = "'73.46"
$ws.Range('E22').Value = '  +2.28%  '
$ws.Range('D23').Value = "'3.34"
$ws.Range('E23').Value = '  +15.63%  '
$ws.Range('D24').Value = "'2.36"
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').Value = "'238.85"
$ws.Range('E25').Value = '  +3.35%  '
$ws.Range('D26').Value = "'8.83"
$ws.Range('E26').Value = '  -4.76%  '
$ws.Range('E27').Value = '  -1.54%  '
$ws.Range('D28').Value = "'11.40"
$ws.Range('E28').Value = '  -5.79%  '
$ws.Range('D29').Value = "'2.19"
$ws.Range('E29').Value = '  -1.42%  '
$ws.Range('D30').Value = "'175.45"
$ws.Range('E30').Value = '  +1.21%  '
$ws.Range('D31').Value = "'36.98"
$ws.Range('E31').Value = '  -9.05%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = "'21.62"
$ws.Range('E32').Value = '  +2.38%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = "'3.12"
$ws.Range('E33').Value = '  -4.45%  '
$ws.Range('D34').Value = "'0.0874"
$ws.Range('E34').Value = '  -3.06%  '
$ws.Range('D35').Value = "'5.65"
$ws.Range('E35').Value = '  +1.41%  '
$ws.Range('E36').Value = '  +5.48%  '
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('E38').Value = '  -2.65%  '
$ws.Range('D39').Value = "'0.0372"
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('E40').Value = '  -2.59%  '
$ws.Range('D41').Value = "'2.39"
$ws.Range('E41').Value = '  -5.43%  '
$ws.Range('D42').Value = "'71.10"
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('D43').Value = "'0.227"
$ws.Range('E43').Value = '  -2.01%  '
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').Value = "'12.21"
$ws.Range('E45').Value = '  -7.71%  '
$ws.Range('E46').Value = '  -2.39%  '
$ws.Range('E47').Value = '  -4.54%  '
$ws.Range('D48').Value = "'1.28"
$ws.Range('E48').Value = '  +2.98%  '
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').Value = "'101.42"
$ws.Range('E50').Value = '  +1.41%  '
$ws.Range('D51').Value = "'1.63"
$ws.Range('E51').Value = '  +5.60%  '
